$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp title (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 07:42"

# Pakistan (row 33) - refreshed case counts
$ws.Range("B33").Value = 1102
$ws.Range("C33").Value = 39
$ws.Range("E33").Value = 1073

# Rows 98-100 were re-sorted (descending by "Casos totales") after the data
# refresh: Kazajistan's updated total (97) now exceeds Camboya's (96) and
# Azerbaiyan's (93), so it moves up to row 98, pushing Camboya to row 99 and
# Azerbaiyan to row 100.

# Row 98: was Camboya -> now Kazajistan (with refreshed stats)
$ws.Range("A98").Value = "Kazajistan"
$ws.Range("B98").Value = 97
$ws.Range("C98").Value = 16
$ws.Range("D98").Value = 2
$ws.Range("E98").Value = 95
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0

# Row 99: was Azerbaiyan -> now Camboya (keeps Camboya's previous stats)
$ws.Range("A99").Value = "Camboya"
$ws.Range("B99").Value = 96
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 10
$ws.Range("E99").Value = 86
$ws.Range("F99").Value = 1
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0

# Row 100: was Kazajistan -> now Azerbaiyan (keeps Azerbaiyan's previous stats)
$ws.Range("A100").Value = "Azerbaiyan"
$ws.Range("B100").Value = 93
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 10
$ws.Range("E100").Value = 81
$ws.Range("F100").Value = 6
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 2
